$d = $word.ActiveDocument

# Locate the short, italic "Zechariah" paragraph that immediately follows the
# "ZEC" heading (Heading2) paragraph, and remove it entirely (including its
# paragraph mark), so the following paragraph (a single space run) takes its
# place right after the "ZEC" heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Zechariah" -and $p.Range.Font.Italic) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the italic 'Zechariah' paragraph to remove."
}

$target.Range.Delete()
